# GIT.docx edit: bold the "instruction" callouts, append "." to the
# first one, and relocate the _GoBack bookmark from the title paragraph
# to the blank paragraph right before "To see the status of a git
# repository".

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the _GoBack bookmark from the title paragraph.
# ---------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------
# 2. "To check the git version" (paragraph 3): bold the existing
#    text, then append a new bold+italic+red "." run.
# ---------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Font.Bold = 1
$insertPoint = $p3.Range
$insertPoint.InsertAfter(".")

# Re-fetch the paragraph range so its End reflects the newly
# inserted character, then locate the "." (the run boundary
# introduces a zero-width marker at the old End, so the actual
# character sits one position earlier).
$p3b = $d.Paragraphs.Item(3)
$dotStart = $p3b.Range.End - 2
$dotRange = $d.Range($dotStart, $dotStart + 1)
$dotRange.Font.Bold = 1

# ---------------------------------------------------------------
# 3. Bold the other single-run instruction paragraphs.
# ---------------------------------------------------------------
$d.Paragraphs.Item(17).Range.Font.Bold = 1
$d.Paragraphs.Item(50).Range.Font.Bold = 1
$d.Paragraphs.Item(54).Range.Font.Bold = 1
$d.Paragraphs.Item(65).Range.Font.Bold = 1

# Paragraph 46 has multiple runs ("To open " / "file explorer where
# this " / "repository" / " is located" / "."); bold the whole range.
$d.Paragraphs.Item(46).Range.Font.Bold = 1

# ---------------------------------------------------------------
# 4. Re-add the _GoBack bookmark to the blank paragraph right
#    before "To see the status of a git repository" (paragraph 64).
# ---------------------------------------------------------------
$p64 = $d.Paragraphs.Item(64)
$d.Bookmarks.Add("_GoBack", $p64.Range)
